$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("googlesearch")

# Update the header text in B1 from "keyword" to "keywords"
$ws.Range("B1").Value = "keywords"

# Activate the sheet and move the selection to B2 (matches authored selection change)
$ws.Activate()
$ws.Range("B2").Select()
